$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.185.59'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '1.632.60'
$ws.Range("E3").Value = '  -0.78%  '
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.52'
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.517'
$ws.Range("E6").Value = '  +1.13%  '
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("E8").Value = '  -0.73%  '
$ws.Range("E9").Value = '  -0.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.26'
$ws.Range("E10").Value = '  +1.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0852'
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("D12").Value = '1.641.55'
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("E14").Value = '  +0.98%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '27.154.31'
$ws.Range("E15").Value = '  -0.08%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.76'
$ws.Range("E16").Value = '  -3.95%  '
$ws.Range("D17").Value = '0.0₃0734'
$ws.Range("E17").Value = '  -0.83%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.95'
$ws.Range("E18").Value = '  -2.05%  '
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.88'
$ws.Range("E20").Value = '  +0.62%  '
$ws.Range("E21").Value = '  -0.93%  '
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("E23").Value = '  -1.57%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '148.25'
$ws.Range("E24").Value = '  +0.14%  '
$ws.Range("E25").Value = '  -0.29%  '
$ws.Range("E26").Value = '  -1.69%  '
$ws.Range("E27").Value = '  -0.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.57'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0504'
$ws.Range("E29").Value = '  -0.56%  '
$ws.Range("E30").Value = '  -0.67%  '
$ws.Range("E31").Value = '  +0.63%  '
$ws.Range("E32").Value = '  -0.88%  '
$ws.Range("D33").Value = '1.318.64'
$ws.Range("E33").Value = '  +4.52%  '
$ws.Range("E34").Value = '  -1.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.44'
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0174'
$ws.Range("E36").Value = '  -1.83%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.846'
$ws.Range("E37").Value = '  +0.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.538'
$ws.Range("E38").Value = '  -0.94%  '
$ws.Range("E39").Value = '  -0.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.26'
$ws.Range("E40").Value = '  +1.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.804'
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '63.83'
$ws.Range("E42").Value = '  +2.70%  '
$ws.Range("D43").Value = '1.770.10'
$ws.Range("E43").Value = '  -1.05%  '
$ws.Range("E44").Value = '  -2.41%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.84'
$ws.Range("E45").Value = '  -1.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.61'
$ws.Range("E46").Value = '  +0.29%  '
$ws.Range("E47").Value = '  +7.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.789'
$ws.Range("E48").Value = '  +16.06%  '
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.56'
$ws.Range("E50").Value = '  -1.78%  '
$ws.Range("E51").Value = '  -0.23%  '
